# Update cryptocurrency price/volume data cells per the latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.441.40"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "3.425.09"
$ws.Range("E3").Value = "  +1.10%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Formula = "'581.14"
$ws.Range("E5").Value = "  -1.38%  "
$ws.Range("D6").Formula = "'177.28"
$ws.Range("E6").Value = "  -2.57%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").Value = "3.419.65"
$ws.Range("E8").Value = "  +1.06%  "
$ws.Range("D9").Formula = "'0.593"
$ws.Range("E9").Value = "  -0.48%  "
$ws.Range("D10").Formula = "'0.198"
$ws.Range("E10").Value = "  +2.07%  "
$ws.Range("E11").Value = "  -0.32%  "
$ws.Range("D12").Formula = "'48.77"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("E13").Value = "  +0.40%  "
$ws.Range("D14").Formula = "'697.89"
$ws.Range("E14").Value = "  +0.58%  "
$ws.Range("D15").Value = "3.977.90"
$ws.Range("E15").Value = "  +1.19%  "
$ws.Range("E16").Value = "  +0.83%  "
$ws.Range("D17").Value = "69.489.22"
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("D18").Value = "3.427.30"
$ws.Range("E18").Value = "  +1.41%  "
$ws.Range("E19").Value = "  +0.83%  "
$ws.Range("D20").Formula = "'17.73"
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("D21").Formula = "'11.40"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Formula = "'0.901"
$ws.Range("E22").Value = "  -0.40%  "
$ws.Range("B23").Value = "Toncoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D23").Formula = "'5.37"
$ws.Range("E23").Value = "  -0.92%  "
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").Formula = "'17.00"
$ws.Range("E24").Value = "  -1.21%  "
$ws.Range("D25").Formula = "'101.35"
$ws.Range("E25").Value = "  -2.95%  "
$ws.Range("E26").Value = "  -1.76%  "
$ws.Range("D27").Formula = "'2.68"
$ws.Range("E27").Value = "  -1.69%  "
$ws.Range("D28").Formula = "'9.65"
$ws.Range("E28").Value = "  +0.17%  "
$ws.Range("D29").Formula = "'33.76"
$ws.Range("E29").Value = "  -0.73%  "
$ws.Range("D30").Formula = "'8.81"
$ws.Range("E30").Value = "  +1.37%  "
$ws.Range("D31").Formula = "'7.01"
$ws.Range("E31").Value = "  -1.31%  "
$ws.Range("D32").Formula = "'3.81"
$ws.Range("E32").Value = "  +3.38%  "
$ws.Range("D33").Formula = "'565.88"
$ws.Range("E33").Value = "  +1.50%  "
$ws.Range("D34").Formula = "'11.04"
$ws.Range("E34").Value = "  -1.26%  "
$ws.Range("E35").Value = "  -1.46%  "
$ws.Range("D36").Formula = "'58.12"
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("D37").Formula = "'0.999"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").Value = "3.618.12"
$ws.Range("E38").Value = "  -3.00%  "
$ws.Range("E39").Value = "  -2.15%  "
$ws.Range("D40").Formula = "'34.99"
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").Value = "0.0₃0739"
$ws.Range("E41").Value = "  +4.51%  "
$ws.Range("D42").Formula = "'3.31"
$ws.Range("E42").Value = "  +2.62%  "
$ws.Range("E43").Value = "  +0.73%  "
$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").Formula = "'3.38"
$ws.Range("E44").Value = "  +3.75%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Formula = "'0.0423"
$ws.Range("E45").Value = "  +0.75%  "
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").Formula = "'0.334"
$ws.Range("E46").Value = "  -1.61%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Formula = "'1.50"
$ws.Range("E47").Value = "  +8.47%  "
$ws.Range("B48").Value = "ThetaToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D48").Formula = "'2.66"
$ws.Range("E48").Value = "  -0.42%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").Formula = "'0.129"
$ws.Range("E49").Value = "  -0.90%  "
$ws.Range("B50").Value = "FirstDigitalUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D50").Formula = "'1.00"
$ws.Range("E50").Value = "  -0.11%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").Formula = "'130.55"
$ws.Range("E51").Value = "  -1.76%  "
